$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CountryRow($row, $name, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($row, 1).Value = $name
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

# --- Block 1: rows 37-38 (Ucrania / Belgica swap with refreshed Belgica data) ---
Set-CountryRow 37 "Belgica" 66662 234 17476 39353 0 11 9833
Set-CountryRow 38 "Ucrania" 66575 0 36744 28202 0 0 1629

# --- Block 2: rows 51-52 (Barein / Honduras swap with refreshed Honduras data) ---
Set-CountryRow 51 "Honduras" 40460 719 5103 34143 0 48 1214
Set-CountryRow 52 "Barein" 39921 0 36531 3249 0 0 141

# --- Block 3: rows 74-75 (El Salvador / Australia swap with refreshed Australia data) ---
Set-CountryRow 74 "Australia" 15580 276 9431 5973 0 9 176
Set-CountryRow 75 "El Salvador" 15446 0 7903 7126 0 0 417

# --- Birmania row 168: refresh a few figures ---
$ws.Cells.Item(168, 2).Value = 351
$ws.Cells.Item(168, 3).Value = 1
$ws.Cells.Item(168, 5).Value = 52

# --- Update the "last refreshed" timestamp banner ---
$ws.Range("A1").Value = "Datos actualizados a 29 de Julio de 2020 a las 05:24"
